# Report utility: add a "FindNewUser" sheet that holds the record-table
# ("Retrieve_Row_Test") data that used to sit at the bottom of AddUser
# (rows 23-25), give it its own header row, and drop those rows from
# AddUser since they now live on their own sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("AddUser")

# New sheet, inserted right after AddUser (i.e. at the end of the tab strip).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws3.Name = "FindNewUser"

# Same header row as the other sheets (Section / Path / Action / Value).
$ws3.Range("A1:D1").Value = $ws2.Range("A1:D1").Value()
$ws3.Range("A1:D1").Font.Bold = $true

# Move the record-table rows (23-25) from AddUser into the new sheet.
$ws3.Range("A2:D4").Value = $ws2.Range("A23:D25").Value()
$ws3.Range("A2:D2").NumberFormat = "@"
$ws3.Range("A3:C3").NumberFormat = "@"

# Match the column widths used elsewhere in the workbook.
$ws3.Columns("A:D").ColumnWidth = 19.14

# Remove the now-duplicated rows from AddUser.
$ws2.Range("A23:D25").EntireRow.Delete()

# Leave the selections/active sheet the way the edited workbook has them.
$ws2.Range("A23:D25").Select()
$ws3.Range("C8").Select()
